$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# --- G1 comment: explain the new SimulationTime triplet format ---
$commentText = "Pavel Balazki:" + "`n" + `
  "Simulation time is defined as time intervals.`n" + `
  "Expected is a triple of values {start, end, resolution}, resolution given in `"points per <time unit>`" as defined in the columne `"SimulationTimeUnit`". Multiple intervals can be separated by a `";`""

$ws.Range("G1").AddComment($commentText) | Out-Null

# --- SimulationTime column values: now triplets <start, end, resolution> instead of a single resolution number ---
$ws.Range("G2").Value = "0, 24, 60"
$ws.Range("G3").Value = "0, 1, 60; 1, 12, 20"
$ws.Range("G4").Value = "0, 12, 20"
$ws.Range("G5").Value = "0, 12, 20"

# --- Restore the active selection on the Scenarios sheet ---
$ws.Activate() | Out-Null
$ws.Range("K22").Select() | Out-Null

Write-Output "Scenarios.xlsx updated: SimulationTime is now a triplet <start, end, resolution>."
